$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is a (cell, new text value) pair taken from the source diff.
# Column D holds price-like strings that are stored as TEXT in the sheet
# (several of them look like plain numbers, e.g. "1.000" or "0.9998"),
# and column E holds percentage strings such as "  -1.48%  ".
# To keep them as text (matching the original cell type) instead of
# letting Excel auto-convert number-looking strings to real numbers,
# we temporarily force a text number format, assign the value, then
# restore the cell to its original "Normal" style/General format so the
# cell's style index is unaffected.
$updates = @(
    @{ Cell = "D2"; Value = '22.066.87' },
    @{ Cell = "E2"; Value = '  -1.48%  ' },
    @{ Cell = "D3"; Value = '1.557.42' },
    @{ Cell = "E3"; Value = '  -0.54%  ' },
    @{ Cell = "D4"; Value = '1.000' },
    @{ Cell = "E4"; Value = '  -0.08%  ' },
    @{ Cell = "D5"; Value = '0.9998' },
    @{ Cell = "E5"; Value = '  -0.10%  ' },
    @{ Cell = "D6"; Value = '287.56' },
    @{ Cell = "E6"; Value = '  +0.51%  ' },
    @{ Cell = "D7"; Value = '0.3874' },
    @{ Cell = "E7"; Value = '  +4.25%  ' },
    @{ Cell = "D8"; Value = '0.3243' },
    @{ Cell = "E8"; Value = '  -0.95%  ' },
    @{ Cell = "D9"; Value = '42.76' },
    @{ Cell = "E9"; Value = '  -7.99%  ' },
    @{ Cell = "D10"; Value = '1.122' },
    @{ Cell = "E10"; Value = '  -1.85%  ' },
    @{ Cell = "D11"; Value = '0.07358' },
    @{ Cell = "E11"; Value = '  -0.60%  ' },
    @{ Cell = "D12"; Value = '1.000' },
    @{ Cell = "E12"; Value = '  -0.07%  ' },
    @{ Cell = "D13"; Value = '19.32' },
    @{ Cell = "E13"; Value = '  -5.26%  ' },
    @{ Cell = "D14"; Value = '5.698' },
    @{ Cell = "E14"; Value = '  -2.36%  ' },
    @{ Cell = "D15"; Value = '6.800' },
    @{ Cell = "E15"; Value = '  -0.23%  ' },
    @{ Cell = "D16"; Value = '0.00001125' },
    @{ Cell = "E16"; Value = '  +2.67%  ' },
    @{ Cell = "D17"; Value = '1.557.08' },
    @{ Cell = "E17"; Value = '  -0.57%  ' },
    @{ Cell = "D18"; Value = '0.06610' },
    @{ Cell = "E18"; Value = '  -1.11%  ' },
    @{ Cell = "D19"; Value = '85.29' },
    @{ Cell = "E19"; Value = '  -0.89%  ' },
    @{ Cell = "D20"; Value = '6.402' },
    @{ Cell = "E20"; Value = '  +1.29%  ' },
    @{ Cell = "E21"; Value = '  -0.06%  ' },
    @{ Cell = "D22"; Value = '15.99' },
    @{ Cell = "E22"; Value = '  -1.44%  ' },
    @{ Cell = "D23"; Value = '11.47' },
    @{ Cell = "E23"; Value = '  -2.42%  ' },
    @{ Cell = "D24"; Value = '22.081.56' },
    @{ Cell = "E24"; Value = '  -1.37%  ' },
    @{ Cell = "D25"; Value = '2.342' },
    @{ Cell = "E25"; Value = '  +1.71%  ' },
    @{ Cell = "D26"; Value = '2.550' },
    @{ Cell = "E26"; Value = '  -0.42%  ' },
    @{ Cell = "D27"; Value = '149.50' },
    @{ Cell = "E27"; Value = '  -0.94%  ' },
    @{ Cell = "D28"; Value = '18.89' },
    @{ Cell = "E28"; Value = '  -1.95%  ' },
    @{ Cell = "D29"; Value = '4.862' },
    @{ Cell = "E29"; Value = '  -1.57%  ' },
    @{ Cell = "D30"; Value = '1.733.06' },
    @{ Cell = "E30"; Value = '  -0.43%  ' },
    @{ Cell = "D31"; Value = '120.72' },
    @{ Cell = "E31"; Value = '  -2.34%  ' },
    @{ Cell = "D32"; Value = '1.111' },
    @{ Cell = "E32"; Value = '  +6.25%  ' },
    @{ Cell = "D33"; Value = '5.829' },
    @{ Cell = "E33"; Value = '  -1.81%  ' },
    @{ Cell = "D34"; Value = '1.705' },
    @{ Cell = "E34"; Value = '  -12.76%  ' },
    @{ Cell = "D35"; Value = '0.08217' },
    @{ Cell = "E35"; Value = '  +0.32%  ' },
    @{ Cell = "D36"; Value = '9.304' },
    @{ Cell = "E36"; Value = '  -3.25%  ' },
    @{ Cell = "D37"; Value = '0.06288' },
    @{ Cell = "E37"; Value = '  +0.07%  ' },
    @{ Cell = "D38"; Value = '0.02303' },
    @{ Cell = "E38"; Value = '  -3.02%  ' },
    @{ Cell = "D39"; Value = '5.230' },
    @{ Cell = "E39"; Value = '  -0.03%  ' },
    @{ Cell = "D40"; Value = '0.2112' },
    @{ Cell = "E40"; Value = '  -3.21%  ' },
    @{ Cell = "D41"; Value = '1.225' },
    @{ Cell = "E41"; Value = '  -7.12%  ' },
    @{ Cell = "D42"; Value = '10.87' },
    @{ Cell = "E42"; Value = '  -2.06%  ' },
    @{ Cell = "D43"; Value = '0.9991' },
    @{ Cell = "E43"; Value = '  -0.15%  ' },
    @{ Cell = "D44"; Value = '0.5963' },
    @{ Cell = "E44"; Value = '  -2.29%  ' },
    @{ Cell = "D45"; Value = '13.54' },
    @{ Cell = "E45"; Value = '  -1.21%  ' },
    @{ Cell = "D46"; Value = '3.717' },
    @{ Cell = "E46"; Value = '  -0.73%  ' },
    @{ Cell = "D47"; Value = '0.5759' },
    @{ Cell = "E47"; Value = '  -3.05%  ' },
    @{ Cell = "E48"; Value = '  -3.51%  ' },
    @{ Cell = "D49"; Value = '119.13' },
    @{ Cell = "E49"; Value = '  -3.54%  ' },
    @{ Cell = "D50"; Value = '1.156' },
    @{ Cell = "E50"; Value = '  -1.69%  ' },
    @{ Cell = "D51"; Value = '0.06902' },
    @{ Cell = "E51"; Value = '  -3.50%  ' }
)

foreach ($u in $updates) {
    $cell = $u.Cell
    $val = $u.Value
    $range = $ws.Range($cell)

    # Detect plain-numeric-looking text (e.g. "1.000", "0.9998") which
    # Excel would otherwise silently coerce into a real number on
    # assignment. Values with thousands separators like "22.066.87" are
    # not valid numbers and are safe to assign directly.
    $looksNumeric = $val -match '^[+-]?\d+(\.\d+)?$'

    if ($looksNumeric) {
        $range.NumberFormat = "@"
        $range.Value = $val
        $range.NumberFormat = "General"
        $range.Style = "Normal"
    } else {
        $range.Value = $val
    }
}
